# Updates column D (Price) and column E (Volume(1h)) cells for the
# cryptos table on the active sheet, per the refreshed GitHub Actions
# scrape. Values are plain text in the workbook (e.g. "27.162.46",
# "  +1.14%  "), so numeric-looking prices are written back through a
# Text-formatted cell (then restored to the Normal style) to keep them
# as text instead of Excel auto-coercing them into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CryptoCellText($cellRef, $text, $forceText) {
    $range = $ws.Range($cellRef)
    if ($forceText) {
        $range.NumberFormat = "@"
        $range.Value = $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

$updates = @(
    @('D2', '27.162.46', $false),
    @('E2', '  +1.14%  ', $false),
    @('D3', '1.893.33', $false),
    @('E3', '  +1.98%  ', $false),
    @('D4', '0.9999', $true),
    @('E4', '  -0.03%  ', $false),
    @('D5', '308.15', $true),
    @('D6', '0.9994', $true),
    @('E6', '  -0.07%  ', $false),
    @('D7', '0.5181', $true),
    @('E7', '  +2.18%  ', $false),
    @('D8', '0.3727', $true),
    @('E8', '  +1.91%  ', $false),
    @('D9', '0.07213', $true),
    @('E9', '  +0.48%  ', $false),
    @('D10', '0.9062', $true),
    @('E10', '  +1.85%  ', $false),
    @('D11', '21.07', $true),
    @('E11', '  +2.10%  ', $false),
    @('D12', '0.07633', $true),
    @('E12', '  +1.37%  ', $false),
    @('D13', '1.896.63', $false),
    @('E13', '  +2.16%  ', $false),
    @('D14', '95.14', $true),
    @('E14', '  +3.76%  ', $false),
    @('D15', '5.294', $true),
    @('E15', '  +1.23%  ', $false),
    @('E16', '  -0.08%  ', $false),
    @('D17', '0.000008521', $true),
    @('E18', '  +2.13%  ', $false),
    @('D19', '0.9995', $true),
    @('E19', '  -0.04%  ', $false),
    @('D20', '27.216.30', $false),
    @('E20', '  +1.15%  ', $false),
    @('D21', '5.062', $true),
    @('E21', '  +0.80%  ', $false),
    @('D22', '2.140.84', $false),
    @('E22', '  +2.68%  ', $false),
    @('D23', '10.64', $true),
    @('E23', '  +3.14%  ', $false),
    @('D24', '6.458', $true),
    @('E24', '  +0.12%  ', $false),
    @('D25', '145.73', $true),
    @('E25', '  -0.25%  ', $false),
    @('D26', '1.792', $true),
    @('E26', '  -0.87%  ', $false),
    @('D27', '18.11', $true),
    @('E27', '  +1.58%  ', $false),
    @('D28', '2.161', $true),
    @('E28', '  +5.32%  ', $false),
    @('D29', '114.74', $true),
    @('E30', '  +5.37%  ', $false),
    @('E31', '  +3.56%  ', $false),
    @('E32', '  +0.08%  ', $false),
    @('D33', '0.05064', $true),
    @('E33', '  -0.62%  ', $false),
    @('E34', '  +4.32%  ', $false),
    @('D35', '0.7633', $true),
    @('E35', '  +4.02%  ', $false),
    @('D36', '3.032', $true),
    @('E36', '  +1.95%  ', $false),
    @('D37', '3.289', $true),
    @('E37', '  +2.27%  ', $false),
    @('D38', '2.570', $true),
    @('E38', '  +3.56%  ', $false),
    @('D39', '0.5641', $true),
    @('E39', '  +6.02%  ', $false),
    @('D40', '0.02002', $true),
    @('E40', '  -0.06%  ', $false),
    @('D41', '1.078', $true),
    @('E41', '  +0.48%  ', $false),
    @('D42', '118.98', $true),
    @('E42', '  +0.37%  ', $false),
    @('D43', '6.605', $true),
    @('E43', '  +1.63%  ', $false),
    @('D44', '8.880', $true),
    @('E44', '  +6.36%  ', $false),
    @('D45', '0.1511', $true),
    @('E45', '  +2.67%  ', $false),
    @('E46', '  +3.76%  ', $false),
    @('D47', '10.16', $true),
    @('E47', '  +2.07%  ', $false),
    @('D48', '0.9992', $true),
    @('E48', '  -0.07%  ', $false),
    @('D49', '1.576', $true),
    @('E49', '  +1.06%  ', $false),
    @('D50', '37.18', $true),
    @('E50', '  +0.57%  ', $false),
    @('D51', '63.66', $true),
    @('E51', '  +1.31%  ', $false)
)

foreach ($u in $updates) {
    Set-CryptoCellText $u[0] $u[1] $u[2]
}
